$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bulk-fill the PSSM data block (B2:K21) with the common background value,
# then overwrite the specific non-background cells individually.
$ws.Range("B2:K21").Value = -17.87045502969748

$ws.Range("C2").Value = 1.885842295917317
$ws.Range("I3").Value = 2.663457341075555
$ws.Range("C4").Value = 2.190569339870655
$ws.Range("D4").Value = 2.17337462599895
$ws.Range("F4").Value = 3.518014814246474
$ws.Range("H4").Value = 1.564973107033493
$ws.Range("J4").Value = 2.09310291354059
$ws.Range("C5").Value = 1.740016558315036
$ws.Range("G5").Value = 2.85874165940897
$ws.Range("B7").Value = 2.615163029505871
$ws.Range("B9").Value = 3.794207954772999
$ws.Range("I10").Value = 1.124287868989161
$ws.Range("K10").Value = 2.037405313791346
$ws.Range("G11").Value = 2.907266626280843
$ws.Range("K11").Value = 1.91751924701138
$ws.Range("E13").Value = 4.321922375414482
$ws.Range("J13").Value = 1.95125823674892
$ws.Range("K13").Value = 1.859207199768961
$ws.Range("D14").Value = 1.185814874062545
$ws.Range("K14").Value = 2.090249757531474
$ws.Range("D15").Value = 0.7245437844447032
$ws.Range("J16").Value = 2.052519829421529
$ws.Range("C17").Value = 2.275388825148554
$ws.Range("D17").Value = 2.73034081673885
$ws.Range("H17").Value = 1.444811750307634
$ws.Range("I17").Value = 2.110912960391523
$ws.Range("J17").Value = 2.185310634908149
$ws.Range("H18").Value = 1.752542290159299
$ws.Range("I18").Value = 1.357614410743433
$ws.Range("J18").Value = 1.664313288202933
$ws.Range("D19").Value = 1.168295205318296
$ws.Range("H19").Value = 1.54761374470555
$ws.Range("I19").Value = 1.499195561323358
$ws.Range("C20").Value = 0.7425034317026615
$ws.Range("D20").Value = 1.421311482246114
$ws.Range("F20").Value = 3.094921325995037
$ws.Range("H20").Value = 1.895001047632825
$ws.Range("I20").Value = 0.82814451810938
$ws.Range("K20").Value = 2.080944952818061
$ws.Range("C21").Value = 0.9142980165266495
$ws.Range("G21").Value = 2.390709056863404
$ws.Range("H21").Value = 2.10664178344667
